$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "splitting an attempt in 2 due to differences":
# row 29 (id 28, DHL, "click to accept a packet") is being split - a new
# sample is appended at the bottom of the table (row 47) that is almost
# identical but carries a more specific description.
$r = 47
$ws.Cells.Item($r, 1).Value = 28
$ws.Cells.Item($r, 2).Value = "msg"

# Match the date formatting already used by the "added" column (style of
# the last data row) rather than introducing a brand new number format.
# (Serial date for 2021-08-08, same as row 46's "added" value.)
$ws.Cells.Item($r, 3).Value = 44416
$ws.Cells.Item(46, 3).Copy()
$ws.Cells.Item($r, 3).PasteSpecial(-4122)

$ws.Cells.Item($r, 4).Value = "MCAST"
$ws.Cells.Item($r, 5).Value = "shortened"
$ws.Cells.Item($r, 6).Value = "delivery"
$ws.Cells.Item($r, 7).Value = "mt"
$ws.Cells.Item($r, 8).Value = "no"
$ws.Cells.Item($r, 9).Value = "click to accept a packet, proper maltese fonts"
$ws.Cells.Item($r, 10).Value = "DHL"

# Selection ends on the cell just below the newly inserted row, matching
# the author's final cursor position after adding the entry.
$ws.Cells.Item($r + 1, 9).Select()
